# Updated symbol list with GitHub Actions - apply cell value corrections
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Set-TextCell "D2" "244.46"
Set-TextCell "D3" "21.95"
Set-TextCell "D4" "5.401"
Set-TextCell "D5" "0.06009"
Set-TextCell "D6" "3.388"
Set-TextCell "D7" "0.8108"
Set-TextCell "D8" "0.9527"
Set-TextCell "D9" "0.1424"
Set-TextCell "D10" "0.07391"
Set-TextCell "D11" "0.03357"
Set-TextCell "D12" "0.03054"
Set-TextCell "D13" "0.09418"
Set-TextCell "D14" "4.004"
Set-TextCell "D15" "0.001602"
Set-TextCell "D16" "0.04830"
Set-TextCell "D17" "0.0005873"
Set-TextCell "D18" "0.006105"
Set-TextCell "D19" "0.005047"
Set-TextCell "D20" "0.0009881"
Set-TextCell "D23" "6.406"
Set-TextCell "D26" "0.1342"
Set-TextCell "B41" "BKEXToken"
Set-TextCell "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D41" "0.1076"
Set-TextCell "E41" "40BKEXTokenBKK"
Set-TextCell "B42" "CEJI"
Set-TextCell "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D42" "0.002720"
Set-TextCell "E42" "41CEJICEJI"
Set-TextCell "B43" "KickToken"
Set-TextCell "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D43" "0.003016"
Set-TextCell "E43" "42KickTokenKICK"
Set-TextCell "D44" "0.005230"
Set-TextCell "D45" "0.00005219"
Set-TextCell "D47" "0.8104"
Set-TextCell "D48" "0.01998"

Write-Output "Applied 37 cell updates"
